# Edit script for reporthouronfcst.xlsx per commit message:
# "alligned hour as index for chart ..> startin from 9 + records of
#  pasthour update only at 00. ELiminated initialization"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Eliminate the initialization row (old row 20 / COV-WELCOME-RES) so the
# table now spans A1:J19 instead of A1:J20.
$ws.Rows("20:20").Delete()

# Refresh the report date shown in the I column header (keep it plain text,
# not an Excel date value, same as the original cell).
$ws.Range("I1").NumberFormat = "@"
$ws.Range("I1").Value2 = "07/04/2023"
$ws.Range("I1").NumberFormat = "General"

# Refresh every queue row (2-19) with the new hourly forecast figures.
# Row 2: COV-GESTIONE-COMM-CMN-RES
$ws.Range("A2").Value = "COV-GESTIONE-COMM-CMN-RES"
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 1.123
$ws.Range("D2").Value = 1.128
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 3
$ws.Range("G2").Value = 18
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 1159
$ws.Range("J2").Value = -99.90267471958585

# Row 3: COV-HELPLINE-CMN-RES
$ws.Range("A3").Value = "COV-HELPLINE-CMN-RES"
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 374
$ws.Range("D3").Value = 378
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 341
$ws.Range("J3").Value = 10.85043988269794

# Row 4: COV-INFOLINE-PRE-RES
$ws.Range("A4").Value = "COV-INFOLINE-PRE-RES"
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 18
$ws.Range("D4").Value = 18
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 15
$ws.Range("J4").Value = 20

# Row 5: COV-INFOPROV-MOBILE-CMN-RES
$ws.Range("A5").Value = "COV-INFOPROV-MOBILE-CMN-RES"
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 243
$ws.Range("D5").Value = 253
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 3
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 294
$ws.Range("J5").Value = -13.94557823129252

# Row 6: COV-INFOPROVISIONING-CMN-RES
$ws.Range("A6").Value = "COV-INFOPROVISIONING-CMN-RES"
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 77
$ws.Range("D6").Value = 78
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 70
$ws.Range("J6").Value = 11.42857142857143

# Row 7: COV-INFOPROVISIONING-RES
$ws.Range("A7").Value = "COV-INFOPROVISIONING-RES"
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 74
$ws.Range("D7").Value = 76
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 111
$ws.Range("J7").Value = -31.53153153153153

# Row 8: COV-MIGRAZIONE-35
$ws.Range("A8").Value = "COV-MIGRAZIONE-35"
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 242
$ws.Range("D8").Value = 249
$ws.Range("E8").Value = 10
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 6
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 407
$ws.Range("J8").Value = -38.82063882063882

# Row 9: COV-MIGRAZIONE-37
$ws.Range("A9").Value = "COV-MIGRAZIONE-37"
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 37
$ws.Range("D9").Value = 36
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 2
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 54
$ws.Range("J9").Value = -33.33333333333334

# Row 10: COV-MOBILE MNP
$ws.Range("A10").Value = "COV-MOBILE MNP"
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 226
$ws.Range("D10").Value = 257
$ws.Range("E10").Value = 34
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 9
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 815
$ws.Range("J10").Value = -68.46625766871166

# Row 11: COV-MOBILE-CMN-RES
$ws.Range("A11").Value = "COV-MOBILE-CMN-RES"
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 274
$ws.Range("D11").Value = 276
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 6
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 376
$ws.Range("J11").Value = -26.59574468085106

# Row 12: COV-MSK-GESTIONE-COMM-RES
$ws.Range("A12").Value = "COV-MSK-GESTIONE-COMM-RES"
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 388
$ws.Range("D12").Value = 578
$ws.Range("E12").Value = 72
$ws.Range("F12").Value = 5
$ws.Range("G12").Value = 6
$ws.Range("H12").Value = 115
$ws.Range("I12").Value = 980
$ws.Range("J12").Value = -41.02040816326531

# Row 13: COV-MSK-HELPLINE-RES
$ws.Range("A13").Value = "COV-MSK-HELPLINE-RES"
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 28
$ws.Range("D13").Value = 28
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 564
$ws.Range("J13").Value = -95.0354609929078

# Row 14: COV-MSK-MOBILE-RES
$ws.Range("A14").Value = "COV-MSK-MOBILE-RES"
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 540
$ws.Range("D14").Value = 658
$ws.Range("E14").Value = 42
$ws.Range("F14").Value = 5
$ws.Range("G14").Value = 10
$ws.Range("H14").Value = 71
$ws.Range("I14").Value = 926
$ws.Range("J14").Value = -28.94168466522678

# Row 15: COV-MSK-RES
$ws.Range("A15").Value = "COV-MSK-RES"
$ws.Range("B15").Value = 0
$ws.Range("C15").Value = 209
$ws.Range("D15").Value = 239
$ws.Range("E15").Value = 30
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 274
$ws.Range("J15").Value = -12.77372262773723

# Row 16: COV-TESTING-RES
$ws.Range("A16").Value = "COV-TESTING-RES"
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 135
$ws.Range("D16").Value = 143
$ws.Range("E16").Value = 10
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 3
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 141
$ws.Range("J16").Value = 1.418439716312059

# Row 17: COV-VIP-RES
$ws.Range("A17").Value = "COV-VIP-RES"
$ws.Range("B17").Value = 0
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 4
$ws.Range("J17").Value = 75

# Row 18: COV-WELCOME-CMN-RES
$ws.Range("A18").Value = "COV-WELCOME-CMN-RES"
$ws.Range("B18").Value = 0
$ws.Range("C18").Value = 19
$ws.Range("D18").Value = 19
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 19
$ws.Range("J18").Value = 0

# Row 19: COV-WELCOME-RES
$ws.Range("A19").Value = "COV-WELCOME-RES"
$ws.Range("B19").Value = 0
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 77
$ws.Range("J19").Value = -83.11688311688312
